$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column S, year 2022 header - same style as the preceding year header (R4)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# New column S, data value for 2022 - based on R5's style, with a new 0.0 number format
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 42
$ws.Range("S5").NumberFormat = "0.0"

$excel.CutCopyMode = 0

# Update selection to match new active cell
$ws.Range("U4").Select()
